# Add speaker notes to slide 2 ("Here we install the HTS-barcode-checker").
#
# Slide 2 currently has no notes slide part at all, so before we can touch
# the notes body text we have to materialize the notes page for that slide
# (PowerPoint normally does this transparently the first time you type into
# the Notes pane). Shapes.AddPlaceholder(2) == ppPlaceholderBody is the
# supported way to get/create the notes body placeholder on a NotesPage; it
# returns the existing placeholder if the slide already has a notes slide
# (harmless no-op), or creates the notesSlide part on first use.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)

$notesPage = $s.NotesPage
$notesBody = $notesPage.Shapes.AddPlaceholder(2)
$notesBody.TextFrame.TextRange.Text = "Here we install the HTS-barcode-checker"
